# Update the "Förändrad" (changed) date column C for rows 2 through 38
# from 2023-10-09 (serial 45208) to 2023-10-13 (serial 45212).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 38; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value2 = 45212
    }
}
